$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text representation
# (values like "26.879.66" or "215.55" must stay as literal text,
# not be auto-converted into numbers) by forcing Text format before
# writing, then restoring the default (Normal) style afterwards so
# the cell styling matches the source workbook.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.879.66"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "1.667.74"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "215.55"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  +4.03%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "0.252"
$ws.Range("E8").Value = "  +1.13%  "
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("D10").Value = "20.20"
$ws.Range("E10").Value = "  +3.36%  "
$ws.Range("D11").Value = "0.0895"
$ws.Range("E11").Value = "  +3.82%  "
$ws.Range("D12").Value = "1.902.84"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").Value = "1.669.52"
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").Value = "66.02"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("D17").Value = "26.888.81"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").Value = "232.01"
$ws.Range("E18").Value = "  -3.73%  "
$ws.Range("D19").Value = "7.81"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("D20").Value = "0.0₃0734"
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "4.45"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "2.21"
$ws.Range("E23").Value = "  -2.65%  "
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("D25").Value = "145.67"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("D33").Value = "1.468.78"
$ws.Range("E33").Value = "  -3.53%  "
$ws.Range("E34").Value = "  +3.37%  "
$ws.Range("D35").Value = "1.63"
$ws.Range("E35").Value = "  +2.81%  "
$ws.Range("D36").Value = "2.42"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "0.573"
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "0.899"
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").Value = "5.82"
$ws.Range("E40").Value = "  -2.29%  "
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").Value = "2.29"
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("E43").Value = "  +6.79%  "
$ws.Range("D44").Value = "65.77"
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("D45").Value = "1.812.44"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("D46").Value = "0.779"
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").Value = "1.53"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("D49").Value = "0.101"
$ws.Range("E49").Value = "  +2.76%  "
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("D51").Value = "7.58"
$ws.Range("E51").Value = "  +0.31%  "

$ws.Range("D2:D51").Style = "Normal"
